$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets (1-indexed, same order as workbook.xml):
#   1 root_hospital_cities
#   2 root_hospital_patients
#   3 root_hospital_users
#   4 attributes
#   5 entities
#   6 packages
# ---------------------------------------------------------------------------
$wsCities   = $wb.Worksheets.Item(1)
$wsPatients = $wb.Worksheets.Item(2)
$wsUsers    = $wb.Worksheets.Item(3)
$wsAttrs    = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------------
# 1) Normalize the example data: display / city names switch from a
#    space-separated form to an underscore-separated "machine" form.
#    ("pape doe" was a typo for papa's display name -> "papa_doe"; the
#    unrelated, correctly spelled "papa doe" on the users sheet is left
#    untouched.)
# ---------------------------------------------------------------------------

# root_hospital_cities: cityName
$wsCities.Cells.Item(2, 1).Value = "new_york"

# root_hospital_patients: displayName (A) and birthplace (E)
$wsPatients.Cells.Item(2, 1).Value = "john_doe"
$wsPatients.Cells.Item(2, 5).Value = "new_york"

$wsPatients.Cells.Item(3, 1).Value = "jane_doe"

$wsPatients.Cells.Item(4, 1).Value = "papa_doe"
$wsPatients.Cells.Item(4, 5).Value = "new_york"

# root_hospital_users: displayName (C)
$wsUsers.Cells.Item(2, 3).Value = "john_doe"
$wsUsers.Cells.Item(3, 3).Value = "jane_doe"
# (row 4 "papa doe" stays as-is)

# ---------------------------------------------------------------------------
# 2) root_hospital_patients: add a new "children" (mref) column between
#    birthplace and disease, and populate the first patient's children.
# ---------------------------------------------------------------------------
$wsPatients.Columns.Item(6).Insert()

$wsPatients.Cells.Item(1, 6).Value = "children"
$wsPatients.Cells.Item(4, 6).Value = "john_doe, jane_doe"

# ---------------------------------------------------------------------------
# 3) attributes: describe the new "children" mref attribute, inserted right
#    after "birthplace" (row 9) and before "disease" (old row 10).
# ---------------------------------------------------------------------------
$wsAttrs.Rows.Item(10).Insert()

$wsAttrs.Cells.Item(10, 1).Value = "children"
$wsAttrs.Cells.Item(10, 2).Value = "root_hospital_patients"
$wsAttrs.Cells.Item(10, 3).Value = "mref"
$wsAttrs.Cells.Item(10, 5).Value = "root_hospital_patients"
$wsAttrs.Cells.Item(10, 7).Value = "children of a patient"

# ---------------------------------------------------------------------------
# 4) Update sheet selections / active cells and move the active tab to the
#    "attributes" sheet.
# ---------------------------------------------------------------------------
$wsCities.Range("E35").Select()
$wsPatients.Range("G15:G16").Select()
$wsUsers.Range("E2").Select()
$wb.Worksheets.Item(5).Range("E21").Select()

$wsAttrs.Activate()
$wsAttrs.Range("E28").Select()
